# 28.01.21 updated ion views - refresh the "Occupancies" analysis sheet with
# the recalculated run: new timestamp, drop the now-unused first data point,
# and overwrite the recalculated occupancy fractions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds the run timestamp used as a simple "last updated" label.
$ws.Range("A1").Value = "28/01/2021 17:03"

# The first occupancy sample (C8) is no longer part of the recalculated
# series, so the cell (value + formatting) is cleared entirely rather than
# just zeroed out.
$ws.Range("C8").Clear()

# Recalculated occupancy fractions for C9:C33 (row -> new value), applied in
# sheet order so the write is deterministic.
$newOccupancies = @(
    @(9,  0.5235014664540886),
    @(10, 0),
    @(11, 0),
    @(12, 0),
    @(13, 0),
    @(14, 0.6711730973214913),
    @(15, 0.682943575243352),
    @(16, 0.4626587186933634),
    @(17, 0.4205932279758092),
    @(18, 0),
    @(19, 0.6344577599912982),
    @(20, 0.782280548001142),
    @(21, 0.8647039252879164),
    @(22, 0.9488675487083401),
    @(23, 0.8698924878045139),
    @(24, 0.9065769118955627),
    @(25, 1),
    @(26, 1),
    @(27, 0.7079042506145522),
    @(28, 1),
    @(29, 1),
    @(30, 1),
    @(31, 0.9525823130225165),
    @(32, 1),
    @(33, 1)
)

foreach ($pair in $newOccupancies) {
    $row = $pair[0]
    $value = $pair[1]
    $ws.Cells.Item($row, 3).Value = $value
}
